# Update countries & provincias Spain
# - Re-sort a few countries whose total-case rank changed (moves their
#   name + figures up/down a row) and refresh the daily case figures
#   for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a full data row (country name + 7 numeric columns) ---
function Set-CountryRow($Row, $Country, $Totales, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Totales
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Iran overtakes Brasil: Iran now row 12 (new figures), Brasil drops to
# row 13 keeping its previous figures.
Set-CountryRow 12 "Iran"   97424 976 78422 12799 2690 47 6203
Set-CountryRow 13 "Brasil" 97100 541 40937 49402 8318 11 6761

# Suiza (row 20): updated totals/new cases/recuperados only.
$ws.Cells.Item(20, 2).Value = 29905
$ws.Cells.Item(20, 3).Value = 88
$ws.Cells.Item(20, 5).Value = 3943

# Rumania (row 37): updated totals/new cases/activos/recuperados/criticos.
$ws.Cells.Item(37, 2).Value = 13163
$ws.Cells.Item(37, 3).Value = 431
$ws.Cells.Item(37, 4).Value = 4869
$ws.Cells.Item(37, 5).Value = 7514
$ws.Cells.Item(37, 6).Value = 255

# Banglades jumps ahead of Dinamarca/Serbia/Filipinas: Banglades gets new
# figures in row 41, while Dinamarca/Serbia/Filipinas each shift down one
# row, carrying their previous figures with them.
Set-CountryRow 41 "Banglades" 9455 665 177  9101 1  2 177
Set-CountryRow 42 "Dinamarca" 9407 0   6889 2043 60 0 475
Set-CountryRow 43 "Serbia"    9362 0   1426 7747 57 0 189
Set-CountryRow 44 "Filipinas" 9223 295 1214 7402 31 4 607

# Marruecos (row 55): updated totals/new cases/activos/recuperados/muertes-hoy/muertes.
$ws.Cells.Item(55, 2).Value = 4880
$ws.Cells.Item(55, 3).Value = 151
$ws.Cells.Item(55, 4).Value = 1424
$ws.Cells.Item(55, 5).Value = 3282
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 174

# Zambia jumps ahead of Togo/Camboya: Zambia gets new figures in row 147,
# while Togo/Camboya each shift down one row, carrying their previous
# figures with them.
Set-CountryRow 147 "Zambia"  124 5 78  43 1 0 3
Set-CountryRow 148 "Togo"    123 0 66  48 0 0 9
Set-CountryRow 149 "Camboya" 122 0 120 2  1 0 0
